$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.489.79'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.28%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.646.43'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.39%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.60'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +6.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.05'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.88%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.608'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +5.84%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.672.77'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.09%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.97%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +5.23%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +7.03%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.123.38'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.480.98'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.36%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.00'
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.673.01'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.98%  '
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000138'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +4.79%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '344.10'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.48'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.57%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.69%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.84'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.28%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.71'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.440'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +5.37%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.995'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.54%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +4.84%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0804'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +9.67%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.71'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.89%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.28'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +6.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '155.92'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.61%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '19.27'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.11'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +5.53%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.907'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +7.66%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +11.96%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +5.61%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.16%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +7.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '303.53'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +7.55%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.86%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.994'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.53%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +4.30%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.27%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +4.24%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.01'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +13.61%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.52'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '10.69'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.04%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +5.08%  '
